# Add 5 new slides to the end of the deck, each a duplicate of the last
# existing slide (Title / Bullets / Subtitle placeholder layout).
$p = $ppt.ActivePresentation

for ($i = 0; $i -lt 5; $i++) {
    $last = $p.Slides.Item($p.Slides.Count)
    $last.Duplicate() | Out-Null
}
